# before config file merging to one file
#
# Merges the separate "DNC1"/"DNC2" duration columns into the main job
# table (rows 3-7), replacing the old start/finish/duration columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): replace start/finish/duration columns with
#     DNC1 / DNC2 / start ---
$ws.Range("M2").Value = "DNC1"
$ws.Range("N2").Value = "DNC2"
$ws.Range("O2").Value = "start"

# --- Row 3 (job 1): M/N become numeric durations, drop old O3 ---
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 40
$ws.Range("O3").Value = ""

# --- Row 4 (job 2): M/N become numeric durations ---
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 40

# --- Row 5 (job 3): fill in the previously-missing D:I columns, add
#     numeric durations ---
$ws.Range("D5").Value = 3001
$ws.Range("E5").Value = 75
$ws.Range("F5").Value = 26
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = "y"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 40

# --- Row 6 (job 4): brand-new row ---
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 4001
$ws.Range("D6").Value = 4001
$ws.Range("E6").Value = 75
$ws.Range("F6").Value = 26
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = "y"
$ws.Range("J6").Value = "AS"
$ws.Range("K6").Value = "Hall"
$ws.Range("L6").Value = "yes"
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 20
$ws.Range("O6").Value = "31.1."

# --- Row 7 (job 5): brand-new row ---
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 5001
$ws.Range("D7").Value = 5001
$ws.Range("E7").Value = 75
$ws.Range("F7").Value = 26
$ws.Range("G7").Value = 38
$ws.Range("H7").Value = 50
$ws.Range("I7").Value = "y"
$ws.Range("J7").Value = "AS"
$ws.Range("K7").Value = "Hall"
$ws.Range("L7").Value = "yes"
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 80

# --- Restore the author's last selection ---
$ws.Range("M11").Select()
